# Edit script for What-can-you-answer.docx
# Applies the commit "[DS] clean up Teacher Resources section, add optional
# activities and links" changes:
#  - "analyze"/"Analyze" -> "relate"/"Relate" (two Q&A sections + two table headers)
#  - italicize+mark "avg" with proofing-error wrapper in 3 question cells
#  - move the stray "_GoBack" bookmark from the very end of the document to
#    just after the second "relate " run

$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 0) Remove the old "_GoBack" bookmark living in the trailing empty
#    paragraph at the end of the document (it gets re-created in step 4
#    right after the second "relate " run). Must run before step 4 inserts
#    the replacement bookmark, because bookmark lookup-by-name is
#    ambiguous while two "_GoBack" bookmarks exist.
# ------------------------------------------------------------------
$oldGoBack = $d.Bookmarks.Item("_GoBack")
$oldGoBack.Delete()

# ------------------------------------------------------------------
# 1) First "analyze " -> "relate" + " " (kept as two runs, both bold)
#    paragraph 55: "...write whether the question is a lookup, compute, or analyze question."
# ------------------------------------------------------------------
$xml1 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
<w:body>
<w:p w14:paraId="3A2C7ABD" w14:textId="77777777" w:rsidR="007A4841" w:rsidRDefault="007A4841" w:rsidP="007A4841"><w:pPr><w:pStyle w:val="Body"/><w:spacing w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:bCs/><w:color w:val="24292E"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:color="24292E"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:b/><w:bCs/><w:color w:val="24292E"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:color="24292E"/></w:rPr><w:t xml:space="preserve">What </w:t></w:r><w:r w:rsidRPr="00BD3E5C"><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:b/><w:bCs/><w:color w:val="24292E"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single" w:color="24292E"/></w:rPr><w:t>can</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:b/><w:bCs/><w:color w:val="24292E"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:color="24292E"/></w:rPr><w:t xml:space="preserve"> you answer? </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:bCs/><w:color w:val="24292E"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:color="24292E"/></w:rPr><w:t xml:space="preserve">For each of the following questions, check the box to the left of questions you </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:bCs/><w:color w:val="24292E"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single" w:color="24292E"/></w:rPr><w:t>can</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:bCs/><w:color w:val="24292E"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:color="24292E"/></w:rPr><w:t xml:space="preserve"> answer. For each </w:t></w:r><w:r w:rsidRPr="00B83593"><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:bCs/><w:i/><w:color w:val="24292E"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:color="24292E"/></w:rPr><w:t>checked</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:bCs/><w:color w:val="24292E"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:color="24292E"/></w:rPr><w:t xml:space="preserve"> question, w</w:t></w:r><w:r w:rsidRPr="00B83593"><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:bCs/><w:color w:val="24292E"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:color="24292E"/></w:rPr><w:t>rite</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:bCs/><w:color w:val="24292E"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:color="24292E"/></w:rPr><w:t xml:space="preserve"> whether the question is a </w:t></w:r><w:r w:rsidRPr="001570FB"><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:b/><w:bCs/><w:color w:val="24292E"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:color="24292E"/></w:rPr><w:t>lookup</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:bCs/><w:color w:val="24292E"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:color="24292E"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:b/><w:bCs/><w:color w:val="24292E"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:color="24292E"/></w:rPr><w:t>compute</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:bCs/><w:color w:val="24292E"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:color="24292E"/></w:rPr><w:t xml:space="preserve">, or </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:b/><w:bCs/><w:color w:val="24292E"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:color="24292E"/></w:rPr><w:t>relate</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:b/><w:bCs/><w:color w:val="24292E"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:color="24292E"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:bCs/><w:color w:val="24292E"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:color="24292E"/></w:rPr><w:t>question.</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$d.Paragraphs.Item(55).Range.InsertXML($xml1)

# ------------------------------------------------------------------
# 2) First "Lookup, Compute or Analyze?" table header -> split into
#    "Lookup, Compute or " + "Relate" + "?"
#    paragraph 59
# ------------------------------------------------------------------
$xml2 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
<w:body>
<w:p w14:paraId="12F6775E" w14:textId="77777777" w:rsidR="007A4841" w:rsidRPr="00BD3E5C" w:rsidRDefault="007A4841" w:rsidP="001E4B26"><w:pPr><w:pStyle w:val="Body"/><w:pBdr><w:top w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:left w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:bottom w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:right w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:between w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:bar w:val="none" w:sz="0" w:color="auto"/></w:pBdr><w:spacing w:line="360" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:b/><w:bCs/><w:color w:val="24292E"/><w:u w:color="24292E"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:b/><w:bCs/><w:color w:val="24292E"/><w:u w:color="24292E"/></w:rPr><w:t xml:space="preserve">Lookup, Compute or </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:b/><w:bCs/><w:color w:val="24292E"/><w:u w:color="24292E"/></w:rPr><w:t>Relate</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:b/><w:bCs/><w:color w:val="24292E"/><w:u w:color="24292E"/></w:rPr><w:t>?</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$d.Paragraphs.Item(59).Range.InsertXML($xml2)

# ------------------------------------------------------------------
# 3) "Create a bar chart showing the avg speed per day" -> italic "avg"
#    wrapped in proofErr spell-check markers
#    paragraph 70
# ------------------------------------------------------------------
$xml3 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
<w:body>
<w:p w14:paraId="11659168" w14:textId="77777777" w:rsidR="007A4841" w:rsidRPr="00BD3E5C" w:rsidRDefault="007A4841" w:rsidP="001E4B26"><w:pPr><w:pStyle w:val="Body"/><w:pBdr><w:top w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:left w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:bottom w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:right w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:between w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:bar w:val="none" w:sz="0" w:color="auto"/></w:pBdr><w:spacing w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:bCs/><w:i/><w:color w:val="24292E"/><w:u w:color="24292E"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:bCs/><w:i/><w:color w:val="24292E"/><w:u w:color="24292E"/></w:rPr><w:t xml:space="preserve">Create a bar chart showing the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:bCs/><w:i/><w:color w:val="24292E"/><w:u w:color="24292E"/></w:rPr><w:t>avg</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:bCs/><w:i/><w:color w:val="24292E"/><w:u w:color="24292E"/></w:rPr><w:t xml:space="preserve"> speed per day</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$d.Paragraphs.Item(70).Range.InsertXML($xml3)

# ------------------------------------------------------------------
# 4) Second "analyze " -> "relate " (single run, text only) + insert the
#    "_GoBack" bookmark right after it
#    paragraph 82
# ------------------------------------------------------------------
$xml4 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
<w:body>
<w:p w14:paraId="5CCAD67D" w14:textId="77777777" w:rsidR="007A4841" w:rsidRPr="00BD3E5C" w:rsidRDefault="007A4841" w:rsidP="007A4841"><w:pPr><w:pStyle w:val="Body"/><w:spacing w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:bCs/><w:i/><w:color w:val="24292E"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single" w:color="24292E"/></w:rPr></w:pPr><w:r w:rsidRPr="00BD3E5C"><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:b/><w:bCs/><w:color w:val="24292E"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:color="24292E"/></w:rPr><w:t xml:space="preserve">What </w:t></w:r><w:r w:rsidRPr="00BD3E5C"><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:b/><w:bCs/><w:i/><w:color w:val="24292E"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single" w:color="24292E"/></w:rPr><w:t>can’t</w:t></w:r><w:r w:rsidRPr="00BD3E5C"><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:b/><w:bCs/><w:color w:val="24292E"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:color="24292E"/></w:rPr><w:t xml:space="preserve"> you answer?</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:bCs/><w:color w:val="24292E"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:color="24292E"/></w:rPr><w:t xml:space="preserve"> For each of the following questions, check the box to the left of questions you </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:bCs/><w:color w:val="24292E"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single" w:color="24292E"/></w:rPr><w:t>cannot</w:t></w:r><w:r w:rsidRPr="00B83593"><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:b/><w:bCs/><w:color w:val="24292E"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:color="24292E"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:bCs/><w:color w:val="24292E"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:color="24292E"/></w:rPr><w:t xml:space="preserve">answer. For each </w:t></w:r><w:r w:rsidRPr="00B83593"><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:bCs/><w:i/><w:color w:val="24292E"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:color="24292E"/></w:rPr><w:t>un-checked</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:bCs/><w:color w:val="24292E"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:color="24292E"/></w:rPr><w:t xml:space="preserve"> question, w</w:t></w:r><w:r w:rsidRPr="00B83593"><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:bCs/><w:color w:val="24292E"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:color="24292E"/></w:rPr><w:t>rite</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:bCs/><w:color w:val="24292E"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:color="24292E"/></w:rPr><w:t xml:space="preserve"> whether the question is a </w:t></w:r><w:r w:rsidRPr="001570FB"><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:b/><w:bCs/><w:color w:val="24292E"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:color="24292E"/></w:rPr><w:t>lookup</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:bCs/><w:color w:val="24292E"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:color="24292E"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:b/><w:bCs/><w:color w:val="24292E"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:color="24292E"/></w:rPr><w:t>compute</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:bCs/><w:color w:val="24292E"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:color="24292E"/></w:rPr><w:t xml:space="preserve">, or </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:b/><w:bCs/><w:color w:val="24292E"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:color="24292E"/></w:rPr><w:t xml:space="preserve">relate </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:bCs/><w:color w:val="24292E"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:color="24292E"/></w:rPr><w:t>question.</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$d.Paragraphs.Item(82).Range.InsertXML($xml4)

# ------------------------------------------------------------------
# 5) Second "Lookup, Compute or Analyze?" table header -> split into
#    "Lookup, Compute or " + "Relate" + "?"
#    paragraph 86
# ------------------------------------------------------------------
$xml5 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
<w:body>
<w:p w14:paraId="7B052088" w14:textId="77777777" w:rsidR="007A4841" w:rsidRPr="00BD3E5C" w:rsidRDefault="007A4841" w:rsidP="001E4B26"><w:pPr><w:pStyle w:val="Body"/><w:pBdr><w:top w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:left w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:bottom w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:right w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:between w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:bar w:val="none" w:sz="0" w:color="auto"/></w:pBdr><w:spacing w:line="360" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:b/><w:bCs/><w:color w:val="24292E"/><w:u w:color="24292E"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:b/><w:bCs/><w:color w:val="24292E"/><w:u w:color="24292E"/></w:rPr><w:t xml:space="preserve">Lookup, Compute or </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:b/><w:bCs/><w:color w:val="24292E"/><w:u w:color="24292E"/></w:rPr><w:t>Relate</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:b/><w:bCs/><w:color w:val="24292E"/><w:u w:color="24292E"/></w:rPr><w:t>?</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$d.Paragraphs.Item(86).Range.InsertXML($xml5)

# ------------------------------------------------------------------
# 6) "What tire pressure produces the highest avg speed?" -> italic "avg"
#    wrapped in proofErr spell-check markers
#    paragraph 89
# ------------------------------------------------------------------
$xml6 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
<w:body>
<w:p w14:paraId="549FEBE2" w14:textId="77777777" w:rsidR="007A4841" w:rsidRPr="00BD3E5C" w:rsidRDefault="007A4841" w:rsidP="001E4B26"><w:pPr><w:pStyle w:val="Body"/><w:pBdr><w:top w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:left w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:bottom w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:right w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:between w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:bar w:val="none" w:sz="0" w:color="auto"/></w:pBdr><w:spacing w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:bCs/><w:i/><w:color w:val="24292E"/><w:u w:color="24292E"/></w:rPr></w:pPr><w:r w:rsidRPr="00BD3E5C"><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:bCs/><w:i/><w:color w:val="24292E"/><w:u w:color="24292E"/></w:rPr><w:t xml:space="preserve">What tire pressure produces the highest </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:bCs/><w:i/><w:color w:val="24292E"/><w:u w:color="24292E"/></w:rPr><w:t>avg</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:bCs/><w:i/><w:color w:val="24292E"/><w:u w:color="24292E"/></w:rPr><w:t xml:space="preserve"> speed?</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$d.Paragraphs.Item(89).Range.InsertXML($xml6)

# ------------------------------------------------------------------
# 7) "What is the avg time it takes this cyclist to ride 1mi?" -> italic
#    "avg" wrapped in proofErr spell-check markers
#    paragraph 93
# ------------------------------------------------------------------
$xml7 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
<w:body>
<w:p w14:paraId="247C8691" w14:textId="77777777" w:rsidR="007A4841" w:rsidRPr="00BD3E5C" w:rsidRDefault="007A4841" w:rsidP="001E4B26"><w:pPr><w:pStyle w:val="Body"/><w:pBdr><w:top w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:left w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:bottom w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:right w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:between w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:bar w:val="none" w:sz="0" w:color="auto"/></w:pBdr><w:spacing w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:bCs/><w:i/><w:color w:val="24292E"/><w:u w:color="24292E"/></w:rPr></w:pPr><w:r w:rsidRPr="00BD3E5C"><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:bCs/><w:i/><w:color w:val="24292E"/><w:u w:color="24292E"/></w:rPr><w:t xml:space="preserve">What is the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:bCs/><w:i/><w:color w:val="24292E"/><w:u w:color="24292E"/></w:rPr><w:t>avg</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:bCs/><w:i/><w:color w:val="24292E"/><w:u w:color="24292E"/></w:rPr><w:t xml:space="preserve"> time it takes this cyclist to ride 1mi?</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$d.Paragraphs.Item(93).Range.InsertXML($xml7)

Write-Host "Edit complete."
